# The sheet's AutoFilter was filtering column C ("lang_code") down to only
# "eng" rows; everything else was present but hidden. This edit removes the
# filter criteria and permanently deletes the rows that were hidden by it,
# so the remaining (previously-visible) rows collapse upward into a
# contiguous A1:F13 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the hidden/filtered-out rows, working bottom-up so the row numbers
# of the groups still above stay valid as each block is removed.
$ws.Rows("15:16").Delete()
$ws.Rows("6:13").Delete()
$ws.Rows("3:4").Delete()

# The leftover empty placeholder cells in column E (no longer meaningful
# once the hidden rows are gone) are cleared out too.
$ws.Range("E2:E7").ClearContents()

# Drop the autofilter criteria entirely (no more filtering on lang_code) and
# resize the filter range down to the new, smaller used range.
$ws.AutoFilterMode = $false
$ws.Range("A1:J13").AutoFilter()

# The filter database named range should track the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='" + $ws.Name + "'!`$A`$1:`$J`$13"
    }
}

# Move the active selection.
$ws.Range("B4").Select()
